$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap column widths: column A gets the old column B width (15.7109375)
# and column B gets the old column A width (16.42578125). The values below
# are the closest settable ColumnWidth inputs that reproduce those widths.
$ws.Columns.Item(1).ColumnWidth = 14.83
$ws.Columns.Item(2).ColumnWidth = 15.67

$ws.Range("A1").Value = -0.28905232522894408
$ws.Range("B1").Value = 0.28839074427983036
$ws.Range("A2").Value = -0.21940476703057499
$ws.Range("B2").Value = 0.2175548554550879
$ws.Range("A3").Value = -0.10428954387691824
$ws.Range("B3").Value = 0.10405723213216689
$ws.Range("A4").Value = -0.096057232240847057
$ws.Range("B4").Value = 0.095642728751487383
$ws.Range("A5").Value = -0.092642728817270203
$ws.Range("B5").Value = 0.091236981951594132
$ws.Range("A6").Value = -0.021895163609949364
$ws.Range("B6").Value = 0.021759204621673689
$ws.Range("A7").Value = -0.011759204771774279
$ws.Range("B7").Value = 0.011740345244322103
$ws.Range("A8").Value = -0.0017403453958841908
$ws.Range("B8").Value = 0.0017321242726393571
$ws.Range("A9").Value = 0.0002678756529341797
$ws.Range("B9").Value = -0.00027002326057479209
$ws.Range("A10").Value = 0.0022700231859662523
$ws.Range("B10").Value = -0.0022693253011407677
$ws.Range("A11").Value = -0.024391975436556024
$ws.Range("B11").Value = 0.02436609802563261
$ws.Range("A12").Value = -0.020866098116241627
$ws.Range("B12").Value = 0.020670606655807333
$ws.Range("A13").Value = -0.017170606752840989
$ws.Range("B13").Value = 0.017081938209120295
$ws.Range("A14").Value = -0.0090819383515805541
$ws.Range("B14").Value = 0.0090532853230973487
$ws.Range("A15").Value = -0.0080532853984873753
$ws.Range("B15").Value = 0.0080348042376563455
$ws.Range("A16").Value = -0.0060348043237565818
$ws.Range("B16").Value = 0.0060032694834024625
$ws.Range("A17").Value = -0.0040032695709157906
$ws.Range("B17").Value = 0.0039999998927822134
$ws.Range("A18").Value = -0.062985977311640795
$ws.Range("B18").Value = 0.06284180144499274
$ws.Range("A19").Value = -0.012091402451409827
$ws.Range("B19").Value = 0.012016186499252957
$ws.Range("A20").Value = -0.0080161865471435334
$ws.Range("B20").Value = 0.0080056369229382085
$ws.Range("A21").Value = -0.0040056369713292739
$ws.Range("B21").Value = 0.0039999999514481743
$ws.Range("A22").Value = -0.10926531164227526
$ws.Range("B22").Value = 0.10853883506624218
$ws.Range("A23").Value = -0.040506526979043223
$ws.Range("B23").Value = 0.040100360756325948
$ws.Range("A24").Value = -0.02010036099058965
$ws.Range("B24").Value = 0.019999999762687182
$ws.Range("A25").Value = -0.060980535368511113
$ws.Range("B25").Value = 0.060935673354642006
$ws.Range("A26").Value = -0.058435673432754243
$ws.Range("B26").Value = 0.058379866357640253
$ws.Range("A27").Value = -0.055879866438323766
$ws.Range("B27").Value = 0.055558808119975822
$ws.Range("A28").Value = -0.053558808208283182
$ws.Range("B28").Value = 0.053352360235045637
$ws.Range("A29").Value = -0.046352360380783608
$ws.Range("B29").Value = 0.046303240102371745
$ws.Range("A30").Value = -0.021169946952621643
$ws.Range("B30").Value = 0.021022058271720656
$ws.Range("A31").Value = -0.014022058426126804
$ws.Range("B31").Value = 0.014000785737287913
$ws.Range("A32").Value = -0.0040007859215265285
$ws.Range("B32").Value = 0.0039999998740185561
